# Mod 3 Project.pptx edit script
# 1. Slide 1 (title slide): set the (empty) center-title placeholder text and
#    delete the now-unused, empty subtitle placeholder shape.
# 2. Slide 9 (Takeaways): reword the first bullet's lead-in run.

$p = $ppt.ActivePresentation

# --- Slide 1: title text + removal of empty subtitle shape -----------------
$s1 = $p.Slides.Item(1)

$titleShape = $s1.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$titleRange.Text = "U.S. Incarceration, Exoneration, and Politics by the Numbers"
$titleRange.LanguageID = "en-GB"

# The second shape is an empty subtitle placeholder that the author removed.
$s1.Shapes.Item(2).Delete()

# --- Slide 9: reword lead-in of the first "Takeaways" bullet ---------------
$s9 = $p.Slides.Item($p.Slides.Count)
$bodyShape = $s9.Shapes.Item(2)
$firstParagraph = $bodyShape.TextFrame.TextRange.Paragraphs(1, 1)
$firstRun = $firstParagraph.Runs(1, 1)
$firstRun.Text = "The data on incarcerated and supervised people are "
